# Apply "added for reporting 30 days data and healty data":
#  - refresh the existing Tuki / Bambang rows (berat/tinggi measurement values)
#  - insert a brand-new balita ("Balita baru tanggal ini")
#  - push the "Balita User 1" record down
#  - append a brand-new balita ("Freya")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Columns like "nikbalita" (D) and "tanggaltimbang" (I) hold digit-only /
    # date-shaped strings. A plain .Value write would let Excel's normal
    # type-inference turn those into numbers/dates, so force Text first and
    # drop the now-unneeded number format once the literal string is in.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-RowData($r, $posyandu, $idbalita, $nik, $nama, $bblahir, $tinggilahir, $alamat, $ortu, $tanggal, $berat, $tinggi) {
    $ws.Range("A$r").Value = $posyandu
    $ws.Range("B$r").Value = $idbalita
    Set-TextCell "D$r" $nik
    $ws.Range("C$r").Value = $nama
    $ws.Range("E$r").Value = $bblahir
    $ws.Range("F$r").Value = $tinggilahir
    $ws.Range("G$r").Value = $alamat
    $ws.Range("H$r").Value = $ortu
    Set-TextCell "I$r" $tanggal
    $ws.Range("J$r").Value = $berat
    $ws.Range("K$r").Value = $tinggi
}

# Row 9 now holds "Tuki" (moved up one row, figures refreshed)
Set-RowData 9 "POSYANDU X" 12.0 "912841" "Tuki" 23.0 10.0 "Rumah lain" "r" "2023-03-27" 8.0 9.0

# Row 10 now holds "Bambang" (moved up one row, figures refreshed)
Set-RowData 10 "POSYANDU X" 2.0 "234634" "Bambang" 5.0 2.0 "Rumah Rafli" "BUKAN RAFLI SUMPA" "2023-02-04" 213.0 1.5

# Row 11: brand-new balita record, added for the 30-day report
Set-RowData 11 "POSYANDU X" 13.0 "12309162309" "Balita baru tanggal ini" 20.0 81.0 "Rumah lain" "r" "2023-05-09" 3.0 55.0

# Row 12: the former "Balita User 1" row, now pushed down
Set-RowData 12 "POSYANDU X" 11.0 "21" "Balita User 1" 2.700000047683716 47.0 "Rumah User 2" "User 2" "2023-03-26" 5.5 2.0

# Row 13: brand-new balita record ("Freya")
Set-RowData 13 "POSYANDU X" 16.0 "31" "Freya" 5.0 40.0 "Rumah lain" "r" "2023-04-27" 123.0 90.0
